# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") holds the strike count that used to be derived from the old
# "Strike#" source; it has been regenerated against the new "K" source and
# every row's value (except row 22, which already matched) needs to be
# rewritten with the freshly calculated number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2;  3  = 1;  4  = 1;  5  = 1;  6  = 2;  7  = 2;  8  = 1;  9  = 0;  10 = 2;
    11 = 1;  12 = 0;  13 = 0;  14 = 0;  15 = 0;  16 = 1;  17 = 0;  18 = 2;  19 = 4;
    20 = 0;  21 = 1;  22 = 0;  23 = 2;  24 = 1;  25 = 2;  26 = 1;  27 = 1;  28 = 3;
    29 = 0;  30 = 0;  31 = 4;  32 = 0;  33 = 0;  34 = 1;  35 = 0;  36 = 2;  37 = 1;
    38 = 2;  39 = 2;  40 = 1;  41 = 1;  42 = 0;  43 = 1;  44 = 1;  45 = 2;  46 = 0;
    47 = 1;  48 = 1;  49 = 1;  50 = 1;  51 = 2;  52 = 3;  53 = 1;  54 = 1;  55 = 0;
    56 = 1;  57 = 4;  58 = 4;  59 = 2;  60 = 2;  61 = 4;  62 = 1;  63 = 1
}

foreach ($row in 2..63) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
